# WeatherStation: add "давление" (pressure) sheet with normalized pressure calc
$wb = $excel.ActiveWorkbook

# Add the new sheet at the end, after the current last sheet ("mlx")
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "давление"

# Title
$ws.Range("B3").Value = "P=P_{0}e^{-Mgh/RT}"

# Inputs
$ws.Range("B5").Value = "h"
$ws.Range("C5").Value = 165
$ws.Range("C5").Font.ThemeColor = 2
$ws.Range("C5").Interior.ThemeColor = 5

$ws.Range("B6").Value = "P0"
$ws.Range("C6").Value = 760
$ws.Range("C6").Font.Color = 192

$ws.Range("B7").Value = "T"
$ws.Range("C7").Value = 20
$ws.Range("C7").Font.Color = 192

# Intermediate calc
$ws.Range("B8").Value = "'-Mgh/RT"
$ws.Range("C8").Formula = '=-0.029*9.81*$C$5/(8.31*(273.15+$C$7))'

$ws.Range("B9").Value = "e^{-Mgh/RT}"
$ws.Range("C9").Formula = '=EXP(C8)'

$ws.Range("B10").Value = "P"
$ws.Range("C10").Formula = '=C6*C9'
$ws.Range("C10").Font.ThemeColor = 2
$ws.Range("C10").Interior.ThemeColor = 1

# Scratch area
$ws.Range("C13").Formula = '=EXP(-0.029*9.81/(8.31*(273.15+$C$7)))'
$ws.Range("F13").Font.ThemeColor = 2

$ws.Range("C14").Formula = '=C13^C5'

# Make the new sheet the active/selected one, matching the authored workbook
$ws.Range("C9").Select()
$wb.Worksheets.Item("давление").Activate()
